$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full roster data for rows 2..19 (row 1 is the header and is untouched)
$data = @(
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Kyshawn George", "SG,SF", "Washington Wizards"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Ja Morant", "PG", "Memphis Grizzlies")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
